$wb = $excel.ActiveWorkbook
$wsDriver = $wb.Worksheets.Item("Driver")

# Fill in new rows 7-9 with the same formatting/pattern as rows 2-6
# (copy format from row 6, then set the literal values)
$wsDriver.Range("A6:F6").Copy()
$wsDriver.Range("A7:F9").PasteSpecial(-4122)

$wsDriver.Range("A7").Value = 3
$wsDriver.Range("B7").Value = 2
$wsDriver.Range("C7").Value = 0
$wsDriver.Range("D7").Value = 4
$wsDriver.Range("E7").Value = 0
$wsDriver.Range("F7").Value = 14

$wsDriver.Range("A8").Value = 4
$wsDriver.Range("B8").Value = 2
$wsDriver.Range("C8").Value = 0
$wsDriver.Range("D8").Value = 4
$wsDriver.Range("E8").Value = 0
$wsDriver.Range("F8").Value = 14

$wsDriver.Range("A9").Value = 5
$wsDriver.Range("B9").Value = 2
$wsDriver.Range("C9").Value = 0
$wsDriver.Range("D9").Value = 4
$wsDriver.Range("E9").Value = 0
$wsDriver.Range("F9").Value = 14

$wsDriver.Rows.Item(7).RowHeight = 15.75
$wsDriver.Rows.Item(8).RowHeight = 15.75
$wsDriver.Rows.Item(9).RowHeight = 15.75
$excel.CutCopyMode = $false

# Select Driver sheet, putting the active selection on E15,
# which also makes Driver the active/selected tab (and Rider no longer selected)
$wsDriver.Activate()
$wsDriver.Range("E15").Select()
